# Issue #99: Fix the G19 observation and other data inconsistencies (#100)
# This updates the regression-summary table on the "Anthropogenic pollution" sheet
# (rows 2-24, columns B:G) to the corrected values from the model re-run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Anthropogenic pollution")

$ws.Range("B2").Value = 2.6061
$ws.Range("C2").Value = 0.3602
$ws.Range("E2").Value = 13.5462
$ws.Range("F2").Value = 6.6873
$ws.Range("G2").Value = 27.4403
$ws.Range("B3").Value = 0.1333
$ws.Range("C3").Value = 0.0922
$ws.Range("D3").Value = 0.1481
$ws.Range("E3").Value = 1.1426
$ws.Range("F3").Value = 0.9538
$ws.Range("G3").Value = 1.3687
$ws.Range("B4").Value = 0.051
$ws.Range("C4").Value = 0.0868
$ws.Range("D4").Value = 0.5572
$ws.Range("E4").Value = 1.0523
$ws.Range("F4").Value = 0.8876
$ws.Range("G4").Value = 1.2475
$ws.Range("B5").Value = -0.2528
$ws.Range("C5").Value = 0.1693
$ws.Range("D5").Value = 0.1354
$ws.Range("E5").Value = 0.7766
$ws.Range("F5").Value = 0.5573
$ws.Range("G5").Value = 1.0823
$ws.Range("B6").Value = 0.2477
$ws.Range("C6").Value = 0.1551
$ws.Range("D6").Value = 0.1103
$ws.Range("E6").Value = 1.2811
$ws.Range("F6").Value = 0.9453
$ws.Range("G6").Value = 1.7362
$ws.Range("B7").Value = -0.5064
$ws.Range("C7").Value = 0.2024
$ws.Range("D7").Value = 0.0124
$ws.Range("E7").Value = 0.6027
$ws.Range("F7").Value = 0.4053
$ws.Range("G7").Value = 0.8962
$ws.Range("B8").Value = 0.345
$ws.Range("C8").Value = 0.1549
$ws.Range("D8").Value = 0.0259
$ws.Range("E8").Value = 1.4119
$ws.Range("F8").Value = 1.0422
$ws.Range("G8").Value = 1.9128
$ws.Range("B9").Value = -0.3807
$ws.Range("C9").Value = 0.2159
$ws.Range("D9").Value = 0.0778
$ws.Range("E9").Value = 0.6834
$ws.Range("F9").Value = 0.4476
$ws.Range("G9").Value = 1.0433
$ws.Range("B10").Value = 0.2623
$ws.Range("C10").Value = 0.1473
$ws.Range("D10").Value = 0.075
$ws.Range("E10").Value = 1.2999
$ws.Range("F10").Value = 0.9739
$ws.Range("G10").Value = 1.7351
$ws.Range("B11").Value = -0.098
$ws.Range("C11").Value = 0.1677
$ws.Range("D11").Value = 0.559
$ws.Range("E11").Value = 0.9066
$ws.Range("F11").Value = 0.6527
$ws.Range("G11").Value = 1.2595
$ws.Range("B12").Value = -0.2421
$ws.Range("C12").Value = 0.2499
$ws.Range("D12").Value = 0.3326
$ws.Range("E12").Value = 0.785
$ws.Range("F12").Value = 0.481
$ws.Range("G12").Value = 1.281
$ws.Range("B13").Value = -0.4407
$ws.Range("C13").Value = 0.3576
$ws.Range("D13").Value = 0.2177
$ws.Range("E13").Value = 0.6436
$ws.Range("F13").Value = 0.3193
$ws.Range("G13").Value = 1.297
$ws.Range("B14").Value = -0.5253
$ws.Range("C14").Value = 0.378
$ws.Range("D14").Value = 0.1647
$ws.Range("E14").Value = 0.5914
$ws.Range("F14").Value = 0.2819
$ws.Range("G14").Value = 1.2407
$ws.Range("B15").Value = -0.5852
$ws.Range("C15").Value = 0.3738
$ws.Range("D15").Value = 0.1174
$ws.Range("E15").Value = 0.557
$ws.Range("F15").Value = 0.2677
$ws.Range("G15").Value = 1.1588
$ws.Range("B16").Value = -0.5458
$ws.Range("C16").Value = 0.3616
$ws.Range("D16").Value = 0.1312
$ws.Range("E16").Value = 0.5794
$ws.Range("F16").Value = 0.2852
$ws.Range("G16").Value = 1.1769
$ws.Range("B17").Value = -0.4765
$ws.Range("C17").Value = 0.3477
$ws.Range("D17").Value = 0.1705
$ws.Range("E17").Value = 0.6209
$ws.Range("F17").Value = 0.3141
$ws.Range("G17").Value = 1.2274
$ws.Range("B18").Value = -0.4667
$ws.Range("C18").Value = 0.3396
$ws.Range("D18").Value = 0.1693
$ws.Range("E18").Value = 0.6271
$ws.Range("F18").Value = 0.3223
$ws.Range("G18").Value = 1.22
$ws.Range("B19").Value = -0.5035
$ws.Range("C19").Value = 0.3396
$ws.Range("D19").Value = 0.1382
$ws.Range("E19").Value = 0.6044
$ws.Range("F19").Value = 0.3106
$ws.Range("G19").Value = 1.176
$ws.Range("B20").Value = -0.4435
$ws.Range("C20").Value = 0.3478
$ws.Range("D20").Value = 0.2022
$ws.Range("E20").Value = 0.6418
$ws.Range("F20").Value = 0.3246
$ws.Range("G20").Value = 1.2688
$ws.Range("B21").Value = -0.3859
$ws.Range("C21").Value = 0.3503
$ws.Range("D21").Value = 0.2705
$ws.Range("E21").Value = 0.6798
$ws.Range("F21").Value = 0.3422
$ws.Range("G21").Value = 1.3506
$ws.Range("B22").Value = -0.5485
$ws.Range("C22").Value = 0.3539
$ws.Range("D22").Value = 0.1212
$ws.Range("E22").Value = 0.5778
$ws.Range("F22").Value = 0.2888
$ws.Range("G22").Value = 1.1563
$ws.Range("B23").Value = -0.7638
$ws.Range("C23").Value = 0.4418
$ws.Range("D23").Value = 0.0838
$ws.Range("E23").Value = 0.4659
$ws.Range("F23").Value = 0.196
$ws.Range("G23").Value = 1.1075
$ws.Range("B24").Value = -0.6584
$ws.Range("E24").Value = 0.5177
$ws.Range("F24").Value = 0.4613
$ws.Range("G24").Value = 0.581
